$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# Update "About" sheet header/version cell
$aboutSheet.Range("A2").Value = "Version: $newVersion"

# Update "About" sheet recommended citation cell
$aboutSheet.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Kostenko Coal Mine, Kazakhstan, M1434, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# Update the version column (S) for each data row on the "Boundaries and methane sources" sheet
for ($row = 2; $row -le 9; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)  # column S = 19
    if ($cell.Value() -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
